$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 01:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1427739
$ws.Range("C4").Value = 19103
$ws.Range("E4").Value = 1034862
$ws.Range("G4").Value = 1616
$ws.Range("H4").Value = 85041

# Row 9 - Brasil
$ws.Range("B9").Value = 188974
$ws.Range("C9").Value = 11372
$ws.Range("D9").Value = 78424
$ws.Range("E9").Value = 97401
$ws.Range("G9").Value = 745
$ws.Range("H9").Value = 13149

# Row 17 - Canada
$ws.Range("B17").Value = 72278
$ws.Range("C17").Value = 1121
$ws.Range("D17").Value = 35164
$ws.Range("E17").Value = 31812
$ws.Range("G17").Value = 133
$ws.Range("H17").Value = 5302

# Row 41 - Colombia
$ws.Range("B41").Value = 12930
$ws.Range("C41").Value = 658
$ws.Range("D41").Value = 3133
$ws.Range("E41").Value = 9288
$ws.Range("G41").Value = 16
$ws.Range("H41").Value = 509

# Row 51 - Chequia
$ws.Range("B51").Value = 8269
$ws.Range("C51").Value = 71
$ws.Range("D51").Value = 5047
$ws.Range("E51").Value = 2932
$ws.Range("G51").Value = 7
$ws.Range("H51").Value = 290

# Row 52 - Noruega
$ws.Range("B52").Value = 8175
$ws.Range("C52").Value = 18
$ws.Range("E52").Value = 7914

# Row 64 - Nigeria
$ws.Range("B64").Value = 4971
$ws.Range("C64").Value = 184
$ws.Range("D64").Value = 1070
$ws.Range("E64").Value = 3743
